# ---------------------------------------------------------------------------
# Variantfile.xlsx edit
#
# Summary of the change (per commit message / diff):
#  1. "main" sheet: mode -> random, N_random -> 1
#  2. New sheet "par-output" inserted between "par-random" and "variants"
#     listing the available OpenFOAM output probes (ID / OUTPUT_NAME / PATH).
#  3. "variants" sheet content regenerated for the new (random, N=1) run:
#     header loses the "mu" column, single data row with new UIn/p values.
#  4. "par-output" becomes the active sheet/tab (selection C4); "main"'s own
#     selection moves to L8.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "main" parameter sheet
# ---------------------------------------------------------------------------
$mainWs = $wb.Worksheets.Item("main")
$mainWs.Range("B3").Value = "random"
$mainWs.Range("B4").Value = 1
$mainWs.Range("L8").Select()

# ---------------------------------------------------------------------------
# 2. Turn the old "variants" sheet into the new "par-output" sheet (this
#    keeps its original sheetId=5 / position, matching the target layout),
#    then create a brand-new "variants" sheet at the end (sheetId=6).
# ---------------------------------------------------------------------------
$outputWs = $wb.Worksheets.Item("variants")
$outputWs.Cells.Clear()
$outputWs.Name = "par-output"

$variantsWs = $wb.Worksheets.Add($null, $outputWs)
$variantsWs.Name = "variants"

# ---------------------------------------------------------------------------
# 3. Populate "par-output"
# ---------------------------------------------------------------------------
$outputWs.Range("A1").Value = "ID"
$outputWs.Range("B1").Value = "OUTPUT_NAME"
$outputWs.Range("C1").Value = "PATH"

$outputWs.Range("A2").Value = 0
$outputWs.Range("B2").Value = "p"
$outputWs.Range("C2").Value = "/postProcessing/probes/0/p"

$outputWs.Range("A3").Value = 1
$outputWs.Range("B3").Value = "magU"
$outputWs.Range("C3").Value = "/postProcessing/probes/0/mag(U)"

$outputWs.Range("A1:C3").RowHeight = 13.8
$outputWs.Columns.Item(2).ColumnWidth = 21.2
$outputWs.Columns.Item(3).ColumnWidth = 24.35

$outputWs.PageSetup.LeftMargin = 56.7
$outputWs.PageSetup.RightMargin = 56.7
$outputWs.PageSetup.TopMargin = 75.8
$outputWs.PageSetup.BottomMargin = 75.8
$outputWs.PageSetup.HeaderMargin = 56.7
$outputWs.PageSetup.FooterMargin = 56.7

# ---------------------------------------------------------------------------
# 4. Populate "variants" (new random-mode result: single row, mu column gone)
# ---------------------------------------------------------------------------
$variantsWs.Range("A1").Value = "ID"
$variantsWs.Range("B1").Value = "UIn"
$variantsWs.Range("C1").Value = "p"
$variantsWs.Range("A1:C1").Font.Bold = $true
$variantsWs.Range("A1:C1").HorizontalAlignment = -4108
$variantsWs.Range("A1:C1").VerticalAlignment = -4160
$variantsWs.Range("A1:C1").Borders.LineStyle = 1

$variantsWs.Range("A2").Value = 0
$variantsWs.Range("B2").Value = 6.231769768424171
$variantsWs.Range("C2").Value = 16968.22813449903

# ---------------------------------------------------------------------------
# 5. Make "par-output" the active tab with C4 selected (matches the diff's
#    bookViews/activeTab + sheetView/selection changes)
# ---------------------------------------------------------------------------
$outputWs.Activate()
$outputWs.Range("C4").Select()
